# Updated cryptos list - Excel COM-interop edit script
# Applies Price (D) and Volume 1h (E) updates per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is a 'clean' decimal number (e.g. 683.10, 0.0000246).
# Pre-format them as Text so Excel keeps the literal digits/trailing zeros
# instead of auto-converting the assignment to a Number (which would drop
# the trailing zero, e.g. turn "683.10" into 683.1).
$textCells = @("D5", "D6", "D11", "D13", "D14", "D18", "D21", "D22", "D23", "D24", "D25", "D26", "D31", "D32", "D33", "D34", "D42", "D46", "D47", "D48", "D49")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "70.539.02"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "3.815.44"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("D5").Value = "683.10"
$ws.Range("E5").Value = "  +8.77%  "
$ws.Range("D6").Value = "170.31"
$ws.Range("E6").Value = "  +2.94%  "
$ws.Range("D7").Value = "3.812.63"
$ws.Range("E7").Value = "  +1.01%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("D11").Value = "7.17"
$ws.Range("E11").Value = "  +5.86%  "
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").Value = "35.91"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").Value = "4.458.36"
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").Value = "3.814.47"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").Value = "70.647.42"
$ws.Range("D18").Value = "17.71"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").Value = "11.22"
$ws.Range("E21").Value = "  +17.40%  "
$ws.Range("D22").Value = "477.03"
$ws.Range("E22").Value = "  +2.00%  "
$ws.Range("D23").Value = "0.714"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("D24").Value = "83.28"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").Value = "0.0000143"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("D26").Value = "12.28"
$ws.Range("E26").Value = "  +2.12%  "
$ws.Range("E27").Value = "  +3.11%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "3.967.22"
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("D31").Value = "2.93"
$ws.Range("E31").Value = "  +9.86%  "
$ws.Range("D32").Value = "2.29"
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("D33").Value = "7.42"
$ws.Range("E33").Value = "  +3.87%  "
$ws.Range("D34").Value = "29.57"
$ws.Range("E34").Value = "  +2.80%  "
$ws.Range("E35").Value = "  +4.46%  "
$ws.Range("E36").Value = "  +2.21%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "3.771.96"
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("E40").Value = "  +2.28%  "
$ws.Range("D42").Value = "0.963"
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("E44").Value = "  +11.51%  "
$ws.Range("D46").Value = "46.01"
$ws.Range("E46").Value = "  +6.46%  "
$ws.Range("D47").Value = "159.50"
$ws.Range("E47").Value = "  +3.29%  "
$ws.Range("D48").Value = "0.000302"
$ws.Range("E48").Value = "  +11.23%  "
$ws.Range("D49").Value = "48.19"
$ws.Range("E49").Value = "  +3.05%  "
$ws.Range("E50").Value = "  +6.19%  "
$ws.Range("E51").Value = "  +1.60%  "
